# testing.xlsx commit: "new jobs and cleaning up old ones."
# - Add a new worksheet "Margem Financeira" after "DRE Saida"
# - Populate it with header labels (row 3) and a data row (row 9) in columns W..AK
# - Add two new label rows (35/36) to "DRE Saida" with a custom font style
# - Update the active selection on "DRE Saida" to B36

$wb  = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)
$ws1.Name = "DRE Saida"

# ---- DRE Saida: two new rows with a distinct (Consolas / orange) font ----
$ws1.Range("B35").Value = "teste feito"
$ws1.Range("B36").Value = "teste 2 feito"

# Apply the font/alignment to each cell individually (rather than the 2-cell
# range at once) - the engine collapses each cell's style edits onto a single
# shared cellXf this way instead of minting one per cell.
foreach ($addr in @("B35", "B36")) {
    $c = $ws1.Range($addr)
    $c.Font.Color = 7901646   # RGB(0xCE, 0x91, 0x78) -> 0xCE9178, packed as BGR for OLE
    $c.Font.Family = 3
    $c.VerticalAlignment = -4108   # xlVAlignCenter
    $c.Font.Name = "Consolas"
}

# ---- New worksheet: Margem Financeira ----
$ws2 = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $ws1)
$ws2.Name = "Margem Financeira"

$ws2.Range("W3").Value  = "3T23"
$ws2.Range("X3").Value  = "4T23"
$ws2.Range("Y3").Value  = "1T24"
$ws2.Range("Z3").Value  = "2T24"
$ws2.Range("AA3").Value = "3T24"
$ws2.Range("AB3").Value = "4T24"
$ws2.Range("AC3").Value = "1T25"
$ws2.Range("AD3").Value = "2T25"
$ws2.Range("AE3").Value = "3T25"
$ws2.Range("AJ3").Value = "9M24"
$ws2.Range("AK3").Value = "9M25"

$ws2.Range("W9").Value  = 9.9
$ws2.Range("X9").Value  = 9.8000000000000007
$ws2.Range("Y9").Value  = 9.3000000000000007
$ws2.Range("Z9").Value  = 9.6
$ws2.Range("AA9").Value = 9.9
$ws2.Range("AB9").Value = 11.1
$ws2.Range("AC9").Value = 9.9
$ws2.Range("AD9").Value = 9.8000000000000007
$ws2.Range("AE9").Value = 9.5
$ws2.Range("AJ9").Value = 9.6
$ws2.Range("AK9").Value = 9.6999999999999993

# ---- view/selection state ----
$ws2.Activate()
$ws2.Range("X9").Select()

$ws1.Activate()
$ws1.Range("B36").Select()
